# "9th Stab - Cosmetic Changes"
#
# The report grows two more "as-of" date columns (Jun_15 / Jun_17) in front
# of the existing trailing date column (which was Jun_13, in column B).
# Column C (the rating/price-target value column) shifts two slots to the
# right, to column E, and the two freshly inserted columns (C, D) are
# pre-filled with the same "UN" placeholder that already sits in column B
# for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns in front of column C. This shifts the old column C
# (price-target/rating values) two columns over to column E, and leaves the
# two new columns (C, D) blank. Column B (the old "Jun_13" header / the "UN"
# placeholder column) is untouched by the insert itself.
$ws.Columns("C:D").Insert()

# The header row: B1 used to hold "Jun_13" - that now belongs in D1 (the
# header shifts along with the data that is two columns to its right), and
# B1/C1 become the two new, more recent, dates.
$oldHeader = $ws.Range("B1").Value2
$ws.Range("D1").Value2 = $oldHeader
$ws.Range("B1").Value2 = "Jun_17"
$ws.Range("C1").Value2 = "Jun_15"

# Every data row (2-27): fill the two newly inserted columns (C, D) with the
# same "UN" placeholder already used in column B.
for ($r = 2; $r -le 27; $r++) {
    $placeholder = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $placeholder
    $ws.Cells.Item($r, 4).Value2 = $placeholder
}

# Cosmetic: give the two new columns the same fixed width as the rest of
# the date columns (matches the sheet's existing 8-character column C).
$ws.Columns("C").ColumnWidth = 7.1666666666666
$ws.Columns("D").ColumnWidth = 7.1666666666666
$ws.Columns("E").ColumnWidth = 7.1666666666666
